# Refresh the crypto price/volume table (cryptos.xlsx) in place.
# Columns: A=rank(unchanged) B=Coin C=Link D=Price E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value2 = "70.382.83"
$ws.Cells.Item(2, 5).Value2 = "  +0.05%  "
# Row 3
$ws.Cells.Item(3, 4).Value2 = "3.608.21"
$ws.Cells.Item(3, 5).Value2 = "  -0.55%  "
# Row 4
$ws.Cells.Item(4, 5).Value2 = "  +0.06%  "
# Row 5
$ws.Cells.Item(5, 4).Value2 = "'581.59"
$ws.Cells.Item(5, 5).Value2 = "  -1.80%  "
# Row 6
$ws.Cells.Item(6, 4).Value2 = "'190.01"
$ws.Cells.Item(6, 5).Value2 = "  -2.04%  "
# Row 7
$ws.Cells.Item(7, 4).Value2 = "3.605.70"
$ws.Cells.Item(7, 5).Value2 = "  -0.48%  "
# Row 8
$ws.Cells.Item(8, 4).Value2 = "'0.630"
$ws.Cells.Item(8, 5).Value2 = "  -2.28%  "
# Row 9
$ws.Cells.Item(9, 5).Value2 = "  +0.04%  "
# Row 10
$ws.Cells.Item(10, 4).Value2 = "'0.188"
$ws.Cells.Item(10, 5).Value2 = "  +3.75%  "
# Row 11
$ws.Cells.Item(11, 4).Value2 = "'0.659"
$ws.Cells.Item(11, 5).Value2 = "  -1.75%  "
# Row 12
$ws.Cells.Item(12, 4).Value2 = "'56.04"
$ws.Cells.Item(12, 5).Value2 = "  -4.08%  "
# Row 13
$ws.Cells.Item(13, 5).Value2 = "  +7.38%  "
# Row 14
$ws.Cells.Item(14, 5).Value2 = "  -2.45%  "
# Row 15
$ws.Cells.Item(15, 4).Value2 = "4.189.51"
$ws.Cells.Item(15, 5).Value2 = "  -0.50%  "
# Row 16
$ws.Cells.Item(16, 4).Value2 = "'19.81"
$ws.Cells.Item(16, 5).Value2 = "  -0.39%  "
# Row 17
$ws.Cells.Item(17, 4).Value2 = "3.620.94"
$ws.Cells.Item(17, 5).Value2 = "  -0.33%  "
# Row 18
$ws.Cells.Item(18, 4).Value2 = "70.394.84"
$ws.Cells.Item(18, 5).Value2 = "  +0.14%  "
# Row 19
$ws.Cells.Item(19, 4).Value2 = "'12.71"
$ws.Cells.Item(19, 5).Value2 = "  -0.20%  "
# Row 20
$ws.Cells.Item(20, 5).Value2 = "  +0.09%  "
# Row 21
$ws.Cells.Item(21, 4).Value2 = "'1.04"
$ws.Cells.Item(21, 5).Value2 = "  -2.08%  "
# Row 22
$ws.Cells.Item(22, 4).Value2 = "'493.58"
$ws.Cells.Item(22, 5).Value2 = "  +1.03%  "
# Row 23
$ws.Cells.Item(23, 4).Value2 = "'19.60"
$ws.Cells.Item(23, 5).Value2 = "  +2.81%  "
# Row 24
$ws.Cells.Item(24, 4).Value2 = "'4.92"
$ws.Cells.Item(24, 5).Value2 = "  -6.85%  "
# Row 25
$ws.Cells.Item(25, 4).Value2 = "'96.29"
$ws.Cells.Item(25, 5).Value2 = "  +5.22%  "
# Row 26
$ws.Cells.Item(26, 5).Value2 = "  -2.03%  "
# Row 27
$ws.Cells.Item(27, 4).Value2 = "'2.98"
$ws.Cells.Item(27, 5).Value2 = "  -5.95%  "
# Row 28
$ws.Cells.Item(28, 4).Value2 = "'11.09"
$ws.Cells.Item(28, 5).Value2 = "  -3.91%  "
# Row 29
$ws.Cells.Item(29, 4).Value2 = "'9.43"
$ws.Cells.Item(29, 5).Value2 = "  -1.72%  "
# Row 30
$ws.Cells.Item(30, 4).Value2 = "'32.23"
$ws.Cells.Item(30, 5).Value2 = "  -2.18%  "
# Row 31
$ws.Cells.Item(31, 4).Value2 = "'7.61"
$ws.Cells.Item(31, 5).Value2 = "  -3.80%  "
# Row 32
$ws.Cells.Item(32, 4).Value2 = "'12.21"
$ws.Cells.Item(32, 5).Value2 = "  -0.51%  "
# Row 33
$ws.Cells.Item(33, 4).Value2 = "'0.117"
$ws.Cells.Item(33, 5).Value2 = "  -3.20%  "
# Row 34
$ws.Cells.Item(34, 4).Value2 = "'65.79"
$ws.Cells.Item(34, 5).Value2 = "  +0.30%  "
# Row 35
$ws.Cells.Item(35, 4).Value2 = "'573.70"
$ws.Cells.Item(35, 5).Value2 = "  -8.54%  "
# Row 36
$ws.Cells.Item(36, 4).Value2 = "'38.48"
$ws.Cells.Item(36, 5).Value2 = "  -5.80%  "
# Row 37
$ws.Cells.Item(37, 4).Value2 = "0.0₃0815"
$ws.Cells.Item(37, 5).Value2 = "  -1.06%  "
# Row 38
$ws.Cells.Item(38, 5).Value2 = "  +0.19%  "
# Row 39
$ws.Cells.Item(39, 4).Value2 = "'3.38"
$ws.Cells.Item(39, 5).Value2 = "  +16.47%  "
# Row 40
$ws.Cells.Item(40, 5).Value2 = "  -3.60%  "
# Row 41
$ws.Cells.Item(41, 5).Value2 = "  +4.17%  "
# Row 42
$ws.Cells.Item(42, 5).Value2 = "  -0.83%  "
# Row 43
$ws.Cells.Item(43, 5).Value2 = "  -6.25%  "
# Row 44
$ws.Cells.Item(44, 5).Value2 = "  -4.66%  "
# Row 45
$ws.Cells.Item(45, 2).Value2 = "Maker"
$ws.Cells.Item(45, 3).Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(45, 4).Value2 = "3.238.40"
$ws.Cells.Item(45, 5).Value2 = "  -1.74%  "
# Row 46
$ws.Cells.Item(46, 2).Value2 = "VeChain"
$ws.Cells.Item(46, 3).Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(46, 4).Value2 = "'0.0444"
$ws.Cells.Item(46, 5).Value2 = "  -2.30%  "
# Row 47
$ws.Cells.Item(47, 2).Value2 = "THORChain"
$ws.Cells.Item(47, 3).Value2 = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Cells.Item(47, 4).Value2 = "'9.82"
$ws.Cells.Item(47, 5).Value2 = "  +6.65%  "
# Row 48
$ws.Cells.Item(48, 2).Value2 = "ApeXProtocol"
$ws.Cells.Item(48, 3).Value2 = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(48, 4).Value2 = "'3.45"
$ws.Cells.Item(48, 5).Value2 = "  +3.64%  "
# Row 49
$ws.Cells.Item(49, 5).Value2 = "  -0.38%  "
# Row 50
$ws.Cells.Item(50, 5).Value2 = "  +0.01%  "
# Row 51
$ws.Cells.Item(51, 4).Value2 = "'3.20"
$ws.Cells.Item(51, 5).Value2 = "  -3.77%  "
